$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3 was stored as a text phone number; convert it to a numeric value.
$ws.Range("C3").Value = 79177131361

# Append new registered user as row 4.
$ws.Range("A4").Value = 974794263
$ws.Range("B4").Value = "Головач Лена"
$ws.Range("C4").Value = "'+992938636344"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "golovach@lena.ker"
$ws.Range("E4").Value = "Buzurgmehr Abdulloev"
